$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.780.33'
$ws.Range("E2").Value = '  +4.67%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.278.95'
$ws.Range("E3").Value = '  +2.17%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.85'
$ws.Range("E5").Value = '  -0.73%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.627'
$ws.Range("E6").Value = '  +0.45%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.40'
$ws.Range("E7").Value = '  +6.13%  '

# Row 8
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
$ws.Range("E9").Value = '  +4.08%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0950'
$ws.Range("E10").Value = '  +4.96%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.65'
$ws.Range("E11").Value = '  -0.70%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.54'
$ws.Range("E12").Value = '  +17.24%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.103'
$ws.Range("E13").Value = '  -0.41%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.619.99'
$ws.Range("E14").Value = '  +2.14%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.69'
$ws.Range("E15").Value = '  -0.05%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.88'
$ws.Range("E16").Value = '  +4.70%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.813'
$ws.Range("E17").Value = '  +1.23%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.279.28'
$ws.Range("E18").Value = '  +1.27%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.697.58'
$ws.Range("E19").Value = '  +4.48%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0941'
$ws.Range("E20").Value = '  +3.69%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.22'
$ws.Range("E21").Value = '  +0.98%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.14'
$ws.Range("E22").Value = '  +0.50%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.88'
$ws.Range("E23").Value = '  +0.89%  '

# Row 24
$ws.Range("E24").Value = '  -0.01%  '

# Row 26
$ws.Range("E26").Value = '  +1.00%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.82'
$ws.Range("E27").Value = '  +1.39%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '171.23'
$ws.Range("E28").Value = '  +1.10%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.137'
$ws.Range("E29").Value = '  -3.99%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.46'
$ws.Range("E30").Value = '  +2.70%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.44'
$ws.Range("E31").Value = '  +2.32%  '

# Row 32
$ws.Range("E32").Value = '  +3.37%  '

# Row 33
$ws.Range("E33").Value = '  +0.18%  '

# Row 34
$ws.Range("E34").Value = '  +6.74%  '

# Row 35
$ws.Range("E35").Value = '  +1.42%  '

# Row 36
$ws.Range("E36").Value = '  +0.39%  '

# Row 37
$ws.Range("E37").Value = '  +0.55%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.67'
$ws.Range("E38").Value = '  +1.20%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.35'
$ws.Range("E39").Value = '  -1.67%  '

# Row 40
$ws.Range("E40").Value = '  +2.93%  '

# Row 41
$ws.Range("E41").Value = '  -0.03%  '

# Row 42
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.96'
$ws.Range("E42").Value = '  +27.20%  '

# Row 43
$ws.Range("B43").Value = 'FTXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.68'
$ws.Range("E43").Value = '  +6.99%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.48'
$ws.Range("E44").Value = '  -1.57%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000221'
$ws.Range("E45").Value = '  -6.19%  '

# Row 46
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0963'
$ws.Range("E46").Value = '  +0.75%  '

# Row 47
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.21'
$ws.Range("E47").Value = '  -1.22%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '97.50'
$ws.Range("E48").Value = '  -1.15%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.481.30'
$ws.Range("E49").Value = '  +0.72%  '

# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.36'
$ws.Range("E50").Value = '  +4.09%  '

# Row 51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.80'
$ws.Range("E51").Value = '  +1.25%  '
